# Edit script: adds "PO Forecast" sheet, renames header columns on the
# two existing sheets, and populates the forecast data.

$wb = $excel.ActiveWorkbook

# --- 1. Rename header labels on existing sheets -----------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet at the end ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header style (bold / border / centered) from an existing sheet's
# header row, then overwrite the text with the new column names.
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:D1"))
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the date-number-format style from the existing date column down
# the new sheet's "ds" column.
$wsWeekly.Range("A2").Copy($wsForecast.Range("A2:A85"))

# --- 3. Populate the forecast data rows ---------------------------------
$data = @(
    @(44934.99999999999,258,-141.2733100188379,643.7596402255991),
    @(44941.99999999999,257,-129.3992409671104,653.1903511057902),
    @(44948.99999999999,256,-154.0959163436046,649.4009982993331),
    @(44955.99999999999,256,-128.1771805814493,645.3413641750457),
    @(44962.99999999999,255,-144.2319857881742,657.3159060994171),
    @(44969.99999999999,254,-128.4320299616246,633.8516243216941),
    @(44976.99999999999,254,-140.824378755799,621.4756309624922),
    @(44983.99999999999,253,-121.1413891686532,612.0655411342221),
    @(44990.99999999999,252,-146.754873680528,654.2266506955992),
    @(44997.99999999999,251,-116.6835667330928,623.3226713588824),
    @(45004.99999999999,251,-152.179775021609,630.0085668563715),
    @(45011.99999999999,250,-139.1723326625701,617.1738932402744),
    @(45018.99999999999,249,-133.6456741177423,604.1736773531439),
    @(45025.99999999999,248,-116.4696281368922,672.8184947462972),
    @(45032.99999999999,248,-111.844759614971,616.7812209568264),
    @(45039.99999999999,247,-124.3620912040832,624.5720646877467),
    @(45046.99999999999,246,-133.0015210765241,625.0216184697824),
    @(45053.99999999999,246,-156.2000555524444,634.1282004373381),
    @(45060.99999999999,245,-133.0874775108635,668.1732803502383),
    @(45067.99999999999,244,-134.227990781636,652.8477474022251),
    @(45074.99999999999,243,-144.2026697526596,624.9058411934104),
    @(45081.99999999999,243,-139.3497673693669,595.2891585122692),
    @(45088.99999999999,242,-127.7701672938973,663.6873662946396),
    @(45095.99999999999,241,-152.8794494678158,663.5368010774176),
    @(45102.99999999999,241,-143.6415588418693,614.0367729853174),
    @(45109.99999999999,240,-151.3926897270968,622.0792011255543),
    @(45116.99999999999,239,-171.9742258984113,612.057420501157),
    @(45123.99999999999,238,-149.0507442087811,630.606445275443),
    @(45130.99999999999,238,-159.2944666043635,602.1629027503682),
    @(45137.99999999999,237,-167.8833381959775,604.8131313288343),
    @(45144.99999999999,236,-151.3866990272836,615.9592020960455),
    @(45151.99999999999,236,-139.2561419987229,635.5931101518844),
    @(45158.99999999999,235,-175.7079736733303,599.5117433094774),
    @(45165.99999999999,234,-174.7617602469248,611.977394978201),
    @(45172.99999999999,233,-150.9215337620121,624.8880251934894),
    @(45179.99999999999,233,-187.8278358079361,598.150341999557),
    @(45186.99999999999,232,-172.8483943184265,591.2627910850578),
    @(45193.99999999999,231,-137.6092477082592,604.3539759144896),
    @(45200.99999999999,230,-172.6960264598256,603.1531585525756),
    @(45207.99999999999,230,-142.8010511987051,584.2697791060189),
    @(45214.99999999999,229,-179.5665528979866,592.0021932210623),
    @(45242.99999999999,226,-173.6563258713912,603.2723391832176),
    @(45249.99999999999,225,-144.4374694214458,613.9131123142288),
    @(45256.99999999999,225,-153.6629634321945,592.1583897444557),
    @(45263.99999999999,224,-182.0023109820127,626.4629243331676),
    @(45270.99999999999,223,-155.3738821796161,614.9630101034097),
    @(45277.99999999999,223,-182.4702311239178,597.3416447168639),
    @(45298.99999999999,220,-193.1688545508984,628.5686477503874),
    @(45305.99999999999,220,-154.6234061949768,631.5209776892103),
    @(45312.99999999999,219,-173.9585248740285,614.3328307011078),
    @(45319.99999999999,218,-167.449276303621,614.6148685719938),
    @(45326.99999999999,218,-186.4384129539426,623.6376185846042),
    @(45333.99999999999,217,-177.5932693180253,617.5910044252859),
    @(45340.99999999999,216,-160.1942411969284,644.5416318673821),
    @(45347.99999999999,215,-171.2661758836055,598.5694310263375),
    @(45354.99999999999,215,-171.6571739055119,581.3142891246982),
    @(45361.99999999999,214,-159.6281211215833,593.0818336918621),
    @(45368.99999999999,213,-186.924363732055,617.607167466743),
    @(45382.99999999999,212,-145.9463412742597,604.4237628151051),
    @(45389.99999999999,211,-192.050153959798,596.7648319824377),
    @(45396.99999999999,210,-175.755487410955,622.3241580006714),
    @(45403.99999999999,210,-181.2555588034299,571.0556521914307),
    @(45410.99999999999,209,-177.4780128048741,593.7745006049596),
    @(45417.99999999999,208,-186.2079025384102,590.4711025905068),
    @(45424.99999999999,207,-186.2160447457654,590.066700602807),
    @(45431.99999999999,207,-152.8105787958874,595.3574646095049),
    @(45438.99999999999,206,-202.705848166788,590.8273341325585),
    @(45480.99999999999,202,-188.5014055084932,594.3000319769188),
    @(45487.99999999999,201,-187.0170144268575,598.1041937417082),
    @(45494.99999999999,200,-174.2887962202338,607.6909979068258),
    @(45515.99999999999,198,-189.3965870098094,603.8341320280974),
    @(45522.99999999999,197,-207.6744510999705,571.5565944466177),
    @(45529.99999999999,197,-201.2090712451327,576.5408978972998),
    @(45536.99999999999,196,-201.1222934407734,607.6474114534784),
    @(45613.99999999999,188,-173.5239205090033,581.1756905119252),
    @(45627.99999999999,187,-203.8242091919407,590.7559328594231),
    @(45634.99999999999,186,-184.5657105553014,586.551371634633),
    @(45641.99999999999,185,-205.553530024584,560.3320177792416),
    @(45648.99999999999,184,-200.10406194288,563.225779095568),
    @(45655.99999999999,184,-235.371224444548,569.9637368916751),
    @(45662.99999999999,183,-201.5278723305055,593.9709459382697),
    @(45669.99999999999,182,-205.9184762569065,575.8989138492578),
    @(45676.99999999999,182,-186.5875769077257,565.8697161797564),
    @(45683.99999999999,181,-223.891425813398,574.4197220214938)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}
